$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the four new columns (H:K) introduced by the migration ---
# Header row
$ws.Range("H1").Value = "Labor Booking User"
$ws.Range("I1").Value = "SiteID"
$ws.Range("J1").Value = "Location ID"
$ws.Range("K1").Value = "Location Number"

# Data rows (same values repeated on row 2 and row 3)
$ws.Range("H2").Value = "a811K0000004fpN"
$ws.Range("I2").Value = "a7q410000004I1W"
$ws.Range("J2").Value = "a7Z4100000000hb"
$ws.Range("K2").Value = "SY_ReceiptLoc"

$ws.Range("H3").Value = "a811K0000004fpN"
$ws.Range("I3").Value = "a7q410000004I1W"
$ws.Range("J3").Value = "a7Z4100000000hb"
$ws.Range("K3").Value = "SY_ReceiptLoc"

# --- Swap the boolean flag in column G between rows 2 and 3 ---
$ws.Range("G2").Value = $false
$ws.Range("G3").Value = $true

# --- Column widths: reflect the re-fitted widths from the source workbook ---
$ws.Columns.Item(2).ColumnWidth = 20.666666666666668
$ws.Columns.Item(8).ColumnWidth = 16.166666666666668
$ws.Columns.Item(9).ColumnWidth = 15.666666666666666
$ws.Columns.Item(10).ColumnWidth = 15.166666666666666
$ws.Columns.Item(11).ColumnWidth = 14.333333333333332

# --- Update the active selection shown in the sheet view ---
$ws.Range("G7").Select()
